# Edit script for LOT2046.xlsx
# Fills in real Portuguese content for Objetivos / Programa resumido / Programa /
# Metodo / Criterio / Norma de recuperacao / Bibliografia (which previously held
# placeholder / misplaced values), and inserts a new row for
# "Docentes responsaveis:" so its value ("1304060 - Maria das Gracas de Almeida
# Felipe") is no longer crammed into the "Objetivos:" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 13 (pushes "Programa resumido:" and everything
#    below it down by one row) to hold the "Docentes responsaveis:" data.
$ws.Rows.Item(13).Insert()

# Copy the B:C number/alignment formatting from the row that is now 14 down into
# the freshly inserted row 13, then drop column A back to the default style so
# the new row has no bold "label" cell (matches the target layout).
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("A13").Style = "Normal"

# 2) Fill the new row with the "Docentes responsaveis:" value (previously
#    mis-stored in row 10 alongside "Objetivos:").
$ws.Range("B13").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("C13").Value = '1304060 - Maria das Graças de Almeida Felipe'

# 3) Objetivos: replace the misplaced "Docentes" text with the real objectives.
$ws.Range("B10").Value = 'Fornecer os conhecimentos necessários sobre os aspectos fundamentais de Microbiologia e Bioquímica Microbiana e sua importância nos estudos sobre Ecologia dos Microrganismos. Fornecer conhecimentos sobre o papel e utilização dos microrganismos nos processos biológicos de interesse à Engenharia Ambiental.'
$ws.Range("C10").Value = 'Fornecer os conhecimentos necessários sobre os aspectos fundamentais de Microbiologia e Bioquímica Microbiana e sua importância nos estudos sobre Ecologia dos Microrganismos. Fornecer conhecimentos sobre o papel e utilização dos microrganismos nos processos biológicos de interesse à Engenharia Ambiental.'

# 4) Programa resumido: replace placeholder "Semestral" with the real short syllabus.
$ws.Range("B14").Value = 'Diversidade metabólica; cultivo e crescimento microbiano; isolamento microbiano; ecossistemas microbianos; biorremediação e biodeterioração  microbiana; bioindicadores.'
$ws.Range("C14").Value = 'Diversidade metabólica; cultivo e crescimento microbiano; isolamento microbiano; ecossistemas microbianos; biorremediação e biodeterioração  microbiana; bioindicadores.'

# 5) Programa: replace placeholder date with the real full syllabus text.
$ws.Range("B16").Value = 'Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos,proteínas e ácidos nucleicos.–Diversidade metabólica: Micro-organismos autotróficos e heterotróficos; glicólise; fermentações; respiração; via das pentoses-fosfato; fotossíntese. –Cultivo e crescimento microbiano: Nutrição microbiana; meios de cultura; fatores ambientais; reprodução e crescimento; medidas e controle de crescimento microbiano. –Isolamento microbiano: Técnicas e meios de isolamento.–Ecossistemas microbianos: Diversidade microbiana e ciclos biogeoquímicos. –Biorremediação e biodeterioração microbiana: Lixiviação bacteriana de metais; bioacumulação e biotransformação microbiana de metais; biodegradação de materiais lignocelulósicos; biodegradação de hidrocarbonetos; biodeterioração de monumentos históricos. –Bioindicadores: Bioindicadores de qualidade de água, ar e solo.'
$ws.Range("C16").Value = 'Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos,proteínas e ácidos nucleicos.–Diversidade metabólica: Micro-organismos autotróficos e heterotróficos; glicólise; fermentações; respiração; via das pentoses-fosfato; fotossíntese. –Cultivo e crescimento microbiano: Nutrição microbiana; meios de cultura; fatores ambientais; reprodução e crescimento; medidas e controle de crescimento microbiano. –Isolamento microbiano: Técnicas e meios de isolamento.–Ecossistemas microbianos: Diversidade microbiana e ciclos biogeoquímicos. –Biorremediação e biodeterioração microbiana: Lixiviação bacteriana de metais; bioacumulação e biotransformação microbiana de metais; biodegradação de materiais lignocelulósicos; biodegradação de hidrocarbonetos; biodeterioração de monumentos históricos. –Bioindicadores: Bioindicadores de qualidade de água, ar e solo.'

# 6) Metodo: replace placeholder "Docentes" text with the real evaluation method.
$ws.Range("B19").Value = 'Duas provas escritas (P1 e P2) distribuídas no semestre.'
$ws.Range("C19").Value = 'Duas provas escritas (P1 e P2) distribuídas no semestre.'

# 7) Criterio: replace the evaluation-method text with the real grading criterion.
$ws.Range("B20").Value = 'MF=Média finalMF = (P1 + P2) / 2'
$ws.Range("C20").Value = 'MF=Média finalMF = (P1 + P2) / 2'

# 8) Norma de recuperacao: replace the grading-criterion text with the real
#    make-up exam rule.
$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.'

# 9) Bibliografia: replace the make-up exam rule text with the real bibliography.
$ws.Range("B22").Value = 'Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Nelson, D.; Cox, M. Princípios de Bioquímica de Lehninger. Artmed Editora. 6a  Edição, 2014.- Pratt, C.; Cornely, K. Bioquímica essencial. Guanabara Koogan. 1a  Edição, 2006. - Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010.- Cooper, G.M. A Célula – Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edição. 2007.- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. - Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edição. 2010.'
$ws.Range("C22").Value = 'Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Nelson, D.; Cox, M. Princípios de Bioquímica de Lehninger. Artmed Editora. 6a  Edição, 2014.- Pratt, C.; Cornely, K. Bioquímica essencial. Guanabara Koogan. 1a  Edição, 2006. - Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010.- Cooper, G.M. A Célula – Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edição. 2007.- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. - Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edição. 2010.'
